# Update cryptocurrency price (D) and 1h volume change (E) columns
# per the Mon May 15 21:47:57 UTC 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.472.91"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "'1.834.01"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -2.88%  "
$ws.Range("D5").Value = "'315.55"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("D7").Value = "'0.4302"
$ws.Range("E7").Value = "  -2.11%  "
$ws.Range("D8").Value = "'0.3707"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").Value = "'0.07277"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "'0.8674"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("D11").Value = "'21.24"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").Value = "'1.826.61"
$ws.Range("E12").Value = "  -1.79%  "
$ws.Range("D13").Value = "'6.704"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'5.369"
$ws.Range("E14").Value = "  -2.60%  "
$ws.Range("D15").Value = "'0.07088"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").Value = "'87.97"
$ws.Range("E16").Value = "  +3.08%  "
$ws.Range("D17").Value = "'1.007"
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").Value = "'0.000008918"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").Value = "'15.26"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").Value = "'27.469.57"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").Value = "'5.178"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "'10.91"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("D24").Value = "'2.053.52"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").Value = "'2.005"
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").Value = "'153.49"
$ws.Range("E26").Value = "  -3.47%  "
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").Value = "'2.161"
$ws.Range("E28").Value = "  +8.12%  "
$ws.Range("D29").Value = "'5.295"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").Value = "'117.68"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "'0.08856"
$ws.Range("E31").Value = "  -2.50%  "
$ws.Range("D32").Value = "'1.211"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "'0.7687"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").Value = "'4.495"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "'2.905"
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("D36").Value = "'1.005"
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("D37").Value = "'1.123"
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").Value = "'0.01963"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").Value = "'0.05291"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'7.173"
$ws.Range("E40").Value = "  +4.15%  "
$ws.Range("D41").Value = "'2.873"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "'0.1680"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "'0.5093"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D44").Value = "'8.688"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "'10.58"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "'0.4745"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").Value = "'106.38"
$ws.Range("E47").Value = "  -3.77%  "
$ws.Range("D48").Value = "'0.06434"
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("D49").Value = "'1.005"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("D50").Value = "'1.672"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").Value = "'1.830"
$ws.Range("E51").Value = "  -3.38%  "
